$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "30.312.95"
$ws.Range("E2").Value = "  -0.91%  "
Set-TextValue $ws.Range("D3") "1.872.38"
$ws.Range("E3").Value = "  -0.15%  "
Set-TextValue $ws.Range("D5") "243.70"
$ws.Range("E5").Value = "  -1.53%  "
Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  +0.13%  "
Set-TextValue $ws.Range("D7") "0.4689"
$ws.Range("E7").Value = "  -1.04%  "
Set-TextValue $ws.Range("D8") "0.2867"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("E9").Value = "  -0.72%  "
Set-TextValue $ws.Range("D10") "21.86"
$ws.Range("E10").Value = "  -0.57%  "
Set-TextValue $ws.Range("D11") "0.07774"
$ws.Range("E11").Value = "  +0.62%  "
Set-TextValue $ws.Range("D12") "1.877.51"
$ws.Range("E12").Value = "  +0.25%  "
Set-TextValue $ws.Range("D13") "95.05"
$ws.Range("E13").Value = "  -0.96%  "
Set-TextValue $ws.Range("D14") "0.7191"
$ws.Range("E14").Value = "  -3.01%  "
Set-TextValue $ws.Range("D15") "5.119"
$ws.Range("E15").Value = "  -0.94%  "
Set-TextValue $ws.Range("D16") "278.37"
$ws.Range("E16").Value = "  +1.49%  "
Set-TextValue $ws.Range("D17") "30.320.01"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("E18").Value = "  -2.30%  "
Set-TextValue $ws.Range("D19") "1.001"
$ws.Range("E19").Value = "  +0.11%  "
Set-TextValue $ws.Range("D20") "0.000007404"
$ws.Range("E20").Value = "  -0.82%  "
Set-TextValue $ws.Range("D21") "2.129.96"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("E22").Value = "  +0.08%  "
Set-TextValue $ws.Range("D23") "5.208"
$ws.Range("E23").Value = "  +0.27%  "
Set-TextValue $ws.Range("D24") "6.225"
$ws.Range("E24").Value = "  +0.87%  "
Set-TextValue $ws.Range("D25") "163.53"
$ws.Range("E25").Value = "  -0.85%  "
Set-TextValue $ws.Range("D26") "9.004"
$ws.Range("E26").Value = "  -2.14%  "
Set-TextValue $ws.Range("D27") "18.61"
$ws.Range("E27").Value = "  -0.36%  "
Set-TextValue $ws.Range("D28") "1.875"
$ws.Range("E28").Value = "  -1.40%  "
Set-TextValue $ws.Range("D29") "1.330"
$ws.Range("E29").Value = "  -1.25%  "
Set-TextValue $ws.Range("D30") "0.09559"
Set-TextValue $ws.Range("D31") "1.464"
$ws.Range("E31").Value = "  -2.97%  "
Set-TextValue $ws.Range("D32") "4.199"
$ws.Range("E32").Value = "  -0.85%  "
Set-TextValue $ws.Range("D33") "4.078"
$ws.Range("E33").Value = "  -0.11%  "
Set-TextValue $ws.Range("D34") "0.04796"
$ws.Range("E34").Value = "  +0.71%  "
Set-TextValue $ws.Range("D35") "1.117"
$ws.Range("E35").Value = "  -0.21%  "
Set-TextValue $ws.Range("D36") "0.6834"
$ws.Range("E36").Value = "  -1.21%  "
Set-TextValue $ws.Range("D37") "2.708"
$ws.Range("E37").Value = "  -0.38%  "
Set-TextValue $ws.Range("D38") "0.01864"
$ws.Range("E38").Value = "  +1.10%  "
Set-TextValue $ws.Range("D39") "2.810"
$ws.Range("E39").Value = "  +1.92%  "
Set-TextValue $ws.Range("D40") "6.202"
$ws.Range("E40").Value = "  -0.76%  "
Set-TextValue $ws.Range("D41") "74.24"
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D42") "0.4209"
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D43") "1.925"
$ws.Range("E43").Value = "  -1.89%  "
Set-TextValue $ws.Range("D44") "1.000"
$ws.Range("E44").Value = "  +0.02%  "
Set-TextValue $ws.Range("D45") "0.8215"
$ws.Range("E45").Value = "  -1.44%  "
Set-TextValue $ws.Range("D46") "100.51"
$ws.Range("E46").Value = "  -0.84%  "
Set-TextValue $ws.Range("D47") "9.587"
$ws.Range("E47").Value = "  +2.32%  "
Set-TextValue $ws.Range("D48") "34.97"
$ws.Range("E48").Value = "  -0.88%  "
Set-TextValue $ws.Range("D49") "6.901"
$ws.Range("E49").Value = "  -0.85%  "
Set-TextValue $ws.Range("D50") "892.68"
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("E51").Value = "  +0.97%  "
